$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Simulate "read from files to excel": each source text file's content
# (including its trailing newline, except for the very last files which
# had no trailing newline) is written verbatim into the corresponding cell.

$rows = @(
    @{ Dong = 1; Files = @("Dong 1 file 1`n", "Dong 1 file 2`n", "Dong 1 file 3`n") },
    @{ Dong = 2; Files = @("Dong 2 file 1`n", "Dong 2 file 2`n", "Dong 2 file 3`n") },
    @{ Dong = 3; Files = @("Dong 3 file 1`n", "Dong 3 file 2`n", "Dong 3 file 3`n") },
    @{ Dong = 4; Files = @("Dong 4 file 1`n", "Dong 4 file 2`n", "Dong 4 file 3`n") },
    @{ Dong = 5; Files = @("Dong 5 file 1",   "Dong 5 file 2`n", "Dong 5 file 3`n") },
    @{ Dong = 6; Files = @($null,             "Dong 6 file 2`n", "Dong 6 file 3`n") },
    @{ Dong = 7; Files = @($null,             "Dong 7 file 2",   $null) }
)

$cols = @("A", "B", "C")

for ($r = 0; $r -lt $rows.Count; $r++) {
    $rowNum = $r + 1
    $files = $rows[$r].Files
    for ($c = 0; $c -lt $cols.Count; $c++) {
        $content = $files[$c]
        if ($null -ne $content) {
            $ws.Range("$($cols[$c])$rowNum").Value = $content
        }
    }
}

[void]$ws.Cells.Select()
